$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value  = "Bidar"
$ws.Range("G4").Value  = "Shivamogga (Shimoga)"
$ws.Range("G8").Value  = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G13").Value = "Shivamogga (Shimoga)"
$ws.Range("G20").Value = "Bidar"
$ws.Range("G26").Value = "Shivamogga (Shimoga)"
$ws.Range("G29").Value = "Shivamogga (Shimoga)"
$ws.Range("G31").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G32").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G33").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G34").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G37").Value = "Bagalkot"
$ws.Range("G41").Value = "Vijayapura (Bijapur)"
$ws.Range("G44").Value = "Shivamogga (Shimoga)"
$ws.Range("G47").Value = "Chikkamagaluru (Chikmagalur)"
